# datacommentdemo.xlsx -- add Birthday / Sex / Working Time columns to the
# employee table, split the old "$init{...} ${employee.bonus}" cell into a
# plain value + a separate comment, move the sub-total / grand-total
# comments to their new columns, and widen the print area accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for two new columns (Sex, Working Time) right after the
#    existing "Birthday" (old "Age") column. Excel shifts every formula,
#    merged range and style that lives at/after column C two slots to the
#    right as part of the insert.
$ws.Columns("C:D").Insert()

# 2. Re-label the header row.
$ws.Range("B7").Value = "Birthday"
$ws.Range("C7").Value = "Sex"
$ws.Range("D7").Value = "Working Time"

# 3. Re-point the template placeholders on the data row.
$ws.Range("B8").Value = '${employee.birthDate}'
$ws.Range("C8").Value = '${employee.sex}'
$ws.Range("D8").Value = '${employee.worktime}'

# The old D8 cell (now shifted to F8) used to hold
# "$init{employee.initbonus} ${employee.bonus}" in a single cell; split it
# into a plain value plus a comment carrying the $init{} directive.
$ws.Range("F8").Value = '${employee.bonus}'

# 4. Move the comments that lived on the old C9/E9/E8 cells to their new
#    homes (column inserts do not relocate existing comments), and add the
#    brand-new widget/init comments the diff introduces.
$ws.Range("E8").Comment.Delete()
$ws.Range("C9").Comment.Delete()
$ws.Range("E9").Comment.Delete()

$ws.Range("B8").AddComment('$widget.calendar{showOn="button" pattern="yyyy/MM/dd" readonlyInput="true"}')
$ws.Range("C8").AddComment('$widget.dropdown{itemLabels="Male;Female" itemValues="M;F" }')
$ws.Range("D8").AddComment('$widget.inputnumber{symbol=" years" symbolPosition="s" minValue="0" maxValue="999" decimalPlaces="2"}')
$ws.Range("F8").AddComment('$init{employee.initbonus}')
$ws.Range("G8").AddComment(' $save{employee.total}')
$ws.Range("E9").AddComment('${department.paymentsum}')
$ws.Range("G9").AddComment('${department.totalsum}')

# 5. Column widths for the two newly inserted columns (closest values this
#    engine's character-width rounding can reach to the authored widths of
#    9.6328125 / 14.36328125).
$ws.Columns("C").ColumnWidth = 8.9
$ws.Columns("D").ColumnWidth = 13.65

# 6. The header banner cells used to merge B:D; now that the table is two
#    columns wider it needs to merge B:G instead of the insert's natural
#    B:F shift.
$ws.Range("B2:F2").UnMerge()
$ws.Range("B2:G2").Merge()
$ws.Range("B3:F3").UnMerge()
$ws.Range("B3:G3").Merge()
$ws.Range("B4:F4").UnMerge()
$ws.Range("B4:G4").Merge()
$ws.Range("B5:F5").UnMerge()
$ws.Range("B5:G5").Merge()

# 7. The print area now spans two more columns (F -> H).
$wb.Names.Item("My_Print_Area").RefersTo = "=Template!`$A`$1:`$H`$10"

# 8. The active selection moved to the first newly-templated cell.
[void]$ws.Range("B8").Select()
